$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 295
$ws.Range("F5").Value = 2066
$ws.Range("F6").Value = 86
$ws.Range("F7").Value = 581
$ws.Range("F8").Value = 446
$ws.Range("F9").Value = 230
$ws.Range("F10").Value = 7751
$ws.Range("F11").Value = 214
$ws.Range("F12").Value = 578
$ws.Range("F13").Value = 822
$ws.Range("F14").Value = 8
$ws.Range("F15").Value = 1894
$ws.Range("F16").Value = 80
$ws.Range("F17").Value = 3244
$ws.Range("F18").Value = 167
$ws.Range("F19").Value = 26
$ws.Range("F20").Value = 61
$ws.Range("F21").Value = 118
$ws.Range("F22").Value = 190
$ws.Range("F23").Value = 137
$ws.Range("F24").Value = 63
$ws.Range("F25").Value = 199
$ws.Range("F26").Value = 86
$ws.Range("F27").Value = 992
$ws.Range("F28").Value = 249
$ws.Range("F29").Value = 4155

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 17

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 139
$ws.Range("F3").Value = 765

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 139
$ws.Range("F4").Value = 765
$ws.Range("F7").Value = 295
$ws.Range("F8").Value = 2067
$ws.Range("F10").Value = 17
$ws.Range("F11").Value = 86
$ws.Range("F12").Value = 581
$ws.Range("F13").Value = 446
$ws.Range("F14").Value = 230
$ws.Range("F15").Value = 7751
$ws.Range("F16").Value = 214
$ws.Range("F17").Value = 578
$ws.Range("F18").Value = 822
$ws.Range("F19").Value = 8
$ws.Range("F20").Value = 1894
$ws.Range("F21").Value = 80
$ws.Range("F22").Value = 3244
$ws.Range("F23").Value = 167
$ws.Range("F24").Value = 26
$ws.Range("F25").Value = 61
$ws.Range("F26").Value = 118
$ws.Range("F27").Value = 190
$ws.Range("F28").Value = 137
$ws.Range("F29").Value = 63
$ws.Range("F30").Value = 199
$ws.Range("F31").Value = 86
$ws.Range("F32").Value = 992
$ws.Range("F33").Value = 249
$ws.Range("F34").Value = 4155
